$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Forms" header cell A2 becomes bold (matches style used by A1/B1 header row) ---
$ws.Range("A2").Font.Bold = $true

# --- Row 4: FieldOID row gets new rich-text Rules content + updated "?" column ---
$c4 = $ws.Range("C4")
$c4.Value2 = "If Forms are downloaded: FORM_OID `nIf CDEs are downloaded: UserID + ""_CDECART""`n" + `
  "Use the caDSR Short Name of the CDE  If no link to a CDE then this column is  null for this row." + `
  "`nOne row for each question on the form UNLESS the Question is a ""all that apply"" type Question.`n" + `
  "For questions with ""all"" in instruction:`n" + `
  "Create the first field for the question with FieldOID as below, + ""_An"" where n a number starting with ""1"" + _LBLAnn"" where nn = n e.g. ""PT_RACE_CD_A1_LBLA1"" `n" + `
  "Create the field for each value in the Question with FieldOID  = same rule as below + ""_An"" where ""n"" is number starting with 1 for the first value, incremented by 1 for each value. E.g. If there are 5 values, _A1, _A2, _A3, _A4, _A5`n" + `
  "e.g. PT_RACE_CD_A1, PT_RACE_CD_A2, PT_RACE_CD_A3, PT_RACE_CD_A4, PT_RACE_CD_A5`n`n`n" + `
  "FieldOID Must be unique within the all the form.`nLimited to 30 characters.`n`n" + `
  "Use the caDSR CDE Short Name that is linked to the Question. Use the Question-CDE link to retreive the CDE short name. `n`n " + `
  "If no link to a CDE then default is UserName + ""_"" + integer (start with 1). `n`n" + `
  "When Questions are used multiple times on a caDSR form, must create these questions as unique fields by adding a sequence number to the  CDE Short Name e.g. " + [char]8220 + "RACE_01" + [char]8221 + ", to ensure uniqueness in Rave.This ensures uniqueness Field OIDs `n`n" + `
  "Replace spaces and other characters:`n  ~ Only underscore characters are permitted`n  ~ Replace space characters with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "." + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "/" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "(" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + ")" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "-" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + [char]8217 + [char]8221 + " with " + [char]8220 + [char]8221

$bold_red = $c4.Characters(80, 96)
$bold_red.Font.Bold = $true
$bold_red.Font.Color = 255

$ws.Range("D4").Value2 = "What to do if length of CD short name = 30 char?"

# --- Row 8 (new): DataFormat ---
$ws.Rows.Item(8).RowHeight = 58
$ws.Range("B8").Value2 = "DataFormat"
$ws.Range("C8").Value2 = "For Questions with ""all"" in the instruction:`nCreate the first field with null DataFormat `nCreate the field for each Value in the Question  DataFormat=1 "
$ws.Range("D8").Value2 = "VD Format: What are conversion rules?`nAs I see if VD Display Format is null  then used '$'+ value of VAL_DOM_MAX_CHAR. What is DataFormat=1 ?"

# --- Row 9 (new): DataDictionaryName ---
$ws.Rows.Item(9).RowHeight = 409.5
$ws.Range("B9").Value2 = "DataDictionaryName"
$ws.Range("C9").Value2 = "For questions with ""all"" in the instruction, there is no entry in this column. `nUse the caDSR CDE Value Domain Long Name truncated to <32 characters plus VD Public ID details:  `nTruncated VD Long Name+" + [char]8220 + "PID" + [char]8221 + "+VDPublicID_" + [char]8221 + "V" + [char]8221 + "+MajorVersion#_+MinorVersion# + ""F""`nAppend an integer in front of ""F"": `nIf VD name ends in a number, such as ""Score 5"", the number corresponds to to the number of times a data dictionary is reused.  IF the dictionary is used only once, ""0F"". Rave requires a unique name for data dictionaries for the protocol.`nTruncate the VD Long Name to <32, and append the letter ""F"" behind any Data dictionary name that ends with a number; prepend an integer in front of the ""F"" for corresponding to the number of times the dictionary is reused, if tje VD is only in the protocol once, this integer will be ""1"".  CDUS_RACE_COD_PID2453600_V6_0_1F`nReplace spaces and other characters:`n  ~ Only underscore characters are permitted`n  ~ Replace space characters with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "." + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "/" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "(" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + ")" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + "-" + [char]8221 + " with " + [char]8220 + "_" + [char]8221 + "`n  ~ Replace " + [char]8220 + [char]8217 + [char]8221 + " with " + [char]8220 + [char]8221

# --- Row 10 (new): ControlType ---
$ws.Rows.Item(10).RowHeight = 72.5
$ws.Range("B10").Value2 = "ControlType"
$ws.Range("C10").Value2 = "For Questions with ""all"" in the Instuction:`nCreate the first field with ControlType =""Text""`nCreate the field for each Value on the Form  ControlType= ""CheckBox"".`nFor Enumerated Question or Enumerated CDE ControlType = ""DropDownList"""

# --- Row 11 (new): PreText ---
$ws.Rows.Item(11).RowHeight = 87
$ws.Range("B11").Value2 = "PreText"
$ws.Range("C11").Value2 = "For Questions with ""all"" in the Instuction:`nCreate the first field with PreText =Form Question Text`nCreate the field for each Value on the Form  PreText = Value  e.g. ""01""  `nUse the Form Question Text, if downloading CDEs, used CDE Preferred Question Text"

# --- Row 12 (new): DefaultValue ---
$ws.Rows.Item(12).RowHeight = 72.5
$ws.Range("B12").Value2 = "DefaultValue"
$ws.Range("C12").Value2 = "if the CDE Value Domain specifies display format, then this format is imported as the FixedUnit.  If no unit of measure display is specified for the CDE VD, then this is blank.`nFor Questions with ""all"" in the Instuction:`nnull FixedUnit"

# --- Row 13 (new): FixedUnit ---
$ws.Rows.Item(13).RowHeight = 72.5
$ws.Range("B13").Value2 = "FixedUnit"
$ws.Range("C13").Value2 = $ws.Range("C12").Value2

# --- View state: scroll to bottom, select D13 ---
$ws.Application.ActiveWindow.ScrollRow = 11
[void]$ws.Range("D13").Select()
